$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto market data.
# D-column values are written via a NumberFormat="@" (Text) round-trip so that
# purely numeric-looking strings (e.g. "1.004") are kept as text, matching the
# original inline-string cell type, instead of Excel auto-coercing them to numbers.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.422.41"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.805.38"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.49%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.90%  "

$ws.Range("E5").Value = "  -0.63%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "306.71"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4518"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3594"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "46.30"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07078"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.8903"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.64%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07801"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "19.44"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.82%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.806.10"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.71%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "5.289"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.20%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "6.323"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "85.10"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.25%  "

$ws.Range("E18").Value = "  -0.89%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000008474"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.68%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.32%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "26.459.94"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.69%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "14.26"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.967"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.047.82"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.96%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "10.52"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.64%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.960"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.51%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "150.77"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "17.81"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.049"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.57%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "111.86"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.88%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.853"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.09%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.08689"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.107"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.11%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.859"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +17.16%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.452"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.87%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7206"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.104"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.075"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.10%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01941"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.57%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.903"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.05095"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.54%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.5114"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.97%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "6.790"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.19%  "

$ws.Range("E45").Value = "  -3.42%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "8.003"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.21%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.4656"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.78%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.54%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.970"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.42%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "100.31"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.90%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.571"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.20%  "

Write-Host "Updated cryptos list values"
